$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AJ1").Value = 0.84629334150465119
$ws.Range("AP1").Value = 0.92907741448500614
$ws.Range("AD2").Value = 0.59209156189274559
$ws.Range("AN3").Value = 0.96796824006150439
$ws.Range("BA3").Value = 0.75391268153871249
$ws.Range("BO3").Value = 0.93956841485179288
$ws.Range("V4").Value = 0.57068717119643042
$ws.Range("F5").Value = 0.68305045175109813
$ws.Range("AP5").Value = 0.80065975539328971
$ws.Range("BG5").Value = 0.80066693117709908
$ws.Range("BI6").Value = 0.88128451719142542
$ws.Range("AS7").Value = 0.94016429542272961
$ws.Range("BG7").Value = 0.78228649898094393
$ws.Range("A9").Value = 0.95605681099017348
$ws.Range("AG9").Value = 0.88723554492911061
$ws.Range("L10").Value = 0.85616196611452111
$ws.Range("Z10").Value = 0.98923183370076029
$ws.Range("BN10").Value = 0.96171705108290795
$ws.Range("D11").Value = 0.87850911722530656
$ws.Range("AE11").Value = 0.78379978158873098
$ws.Range("I12").Value = 0.88106398881894621
$ws.Range("N12").Value = 0.94731625273412678
$ws.Range("BK12").Value = 0.56678477663975346
$ws.Range("I13").Value = 0.97736491258943881
$ws.Range("R14").Value = 0.95007918060826124
$ws.Range("S15").Value = 0.8260196693077787
$ws.Range("AC15").Value = 0.96835137814688987
$ws.Range("AE15").Value = 0.91613960328508348
$ws.Range("V16").Value = 0.83459500883805116
$ws.Range("AT16").Value = 0.91847448439131774
$ws.Range("J17").Value = 0.57515833213092482
$ws.Range("S18").Value = 0.89736286147933808
$ws.Range("AM19").Value = 0.92321492026043084
$ws.Range("AN19").Value = 0.9636372270849326
$ws.Range("R20").Value = 0.99542873213929561
$ws.Range("AR20").Value = 0.97384189302925317
$ws.Range("BG20").Value = 0.9940917822325549
$ws.Range("N21").Value = 0.9851709430570833
$ws.Range("BJ22").Value = 0.96916293635565343
$ws.Range("C23").Value = 0.85657572586668196
$ws.Range("V23").Value = 0.65532971764400805
$ws.Range("BF23").Value = 0.81818704679050702
$ws.Range("AQ24").Value = 0.69976571909869167
$ws.Range("W25").Value = 0.82515131447028611
$ws.Range("X25").Value = 0.82962100979162456
$ws.Range("AG25").Value = 0.92444826752401066
$ws.Range("S26").Value = 0.96449040824178867
$ws.Range("U26").Value = 0.98305823006270587
$ws.Range("AA26").Value = 0.8919057222861515
$ws.Range("N28").Value = 0.94947068275559032
$ws.Range("AV28").Value = 0.96784219213631351
$ws.Range("AD29").Value = 0.93599047524654955
$ws.Range("AE29").Value = 0.99014400348880582
$ws.Range("AR29").Value = 0.99089016184611389
$ws.Range("BM29").Value = 0.97410008387997615
$ws.Range("F30").Value = 0.8150112717316228
$ws.Range("AZ30").Value = 0.98680635196275746
$ws.Range("AA31").Value = 0.81063351775457493
$ws.Range("AD31").Value = 0.99673888302523461
$ws.Range("AZ32").Value = 0.86128788468119188
$ws.Range("BH32").Value = 0.97479179683243855
$ws.Range("F33").Value = 0.92634681423100873
$ws.Range("AJ34").Value = 0.99265807194982392
$ws.Range("BN34").Value = 0.73570882985725761
$ws.Range("AE35").Value = 0.80130953734613253
$ws.Range("AX35").Value = 0.89775367754060142
$ws.Range("H36").Value = 0.7650618831081839
$ws.Range("AB36").Value = 0.88842632974020097
$ws.Range("BD36").Value = 0.99353179585242213
$ws.Range("D37").Value = 0.899661531606059
$ws.Range("T37").Value = 0.78573140683082876
$ws.Range("AV37").Value = 0.98162952559472139
$ws.Range("BE37").Value = 0.79514766203162313
$ws.Range("G38").Value = 0.83244036482212935
$ws.Range("F39").Value = 0.99993846785807095
$ws.Range("AU39").Value = 0.66242840327444674
$ws.Range("X40").Value = 0.93934591462818284
$ws.Range("AJ40").Value = 0.69896920014730068
$ws.Range("AX40").Value = 0.73994822367470769
$ws.Range("BE41").Value = 0.86551626826008721
$ws.Range("BP41").Value = 0.95301584389616456
$ws.Range("Q42").Value = 0.93662791704771653
$ws.Range("BC42").Value = 0.90817647851763805
$ws.Range("AP43").Value = 0.85002648961834493
$ws.Range("AS43").Value = 0.59047916881659734
$ws.Range("AP44").Value = 0.85965931925813832
$ws.Range("AT44").Value = 0.71396972246907353
$ws.Range("B45").Value = 0.79891685449960226
$ws.Range("K45").Value = 0.75388502005869362
$ws.Range("X45").Value = 0.6741935706811284
$ws.Range("B46").Value = 0.75634126527468126
$ws.Range("M47").Value = 0.81244661634471249
$ws.Range("AF47").Value = 0.93785430459514274
$ws.Range("AV47").Value = 0.84603866053321974
$ws.Range("AL48").Value = 0.70947240214933061
$ws.Range("AO48").Value = 0.87960400921164328
$ws.Range("Z49").Value = 0.61580169646369365
$ws.Range("AR49").Value = 0.80113235045040687
$ws.Range("AZ49").Value = 0.75733614605364308
$ws.Range("BE49").Value = 0.85079417707092087
$ws.Range("E50").Value = 0.85020063625482456
$ws.Range("AZ50").Value = 0.94632043520310072
$ws.Range("AO51").Value = 0.84089785349318946
$ws.Range("P53").Value = 0.67856475519298487
$ws.Range("AL53").Value = 0.95119498709856165
$ws.Range("AT53").Value = 0.91065342918865477
$ws.Range("U54").Value = 0.97616286682207032
$ws.Range("AC54").Value = 0.70188767705240496
$ws.Range("AS54").Value = 0.87581730578773032
$ws.Range("C55").Value = 0.88630552417616293
$ws.Range("B56").Value = 0.8750526916903536
$ws.Range("AU56").Value = 0.7427595445287708
$ws.Range("AM57").Value = 0.99269672224092231
$ws.Range("BC57").Value = 0.62364396716829418
$ws.Range("H58").Value = 0.99947629405340677
$ws.Range("BF59").Value = 0.61625854179573158
$ws.Range("AY60").Value = 0.55012683652632277
$ws.Range("BJ60").Value = 0.76913729579706303
$ws.Range("BN60").Value = 0.96855624429095588
$ws.Range("M61").Value = 0.84644752087620989
$ws.Range("BF61").Value = 0.8843436420121763
$ws.Range("O62").Value = 0.81268069023056178
$ws.Range("BL62").Value = 0.95220872361623743
$ws.Range("H63").Value = 0.94918478635686609
$ws.Range("AI63").Value = 0.85536998904668582
$ws.Range("BJ63").Value = 0.67349789740398491
$ws.Range("B64").Value = 0.9186063459054189
$ws.Range("T64").Value = 0.98006485479231631
$ws.Range("M65").Value = 0.9720812740714424
$ws.Range("BD65").Value = 0.69655516068184331
$ws.Range("G66").Value = 0.92370891964023483
$ws.Range("AX66").Value = 0.96287669757544447
$ws.Range("AZ67").Value = 0.86023223060878484
$ws.Range("A68").Value = 0.69178007606457625
$ws.Range("AL68").Value = 0.81077276855865099
$ws.Range("BF68").Value = 0.78910461430001488
